$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "59.639.96"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.78%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.642.54"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.05%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "537.24"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.57%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "145.99"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +3.48%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.13%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.573"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +1.21%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "6.93"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("E10").Value = "  -0.77%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.338"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  +0.34%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.109.51"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +1.51%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "59.527.48"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.69%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "21.37"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +4.15%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.641.43"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("E18").Value = "  +2.99%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "339.47"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("E22").Value = "  -0.04%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "66.31"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  -0.49%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +1.65%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0750"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -0.08%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.66"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -2.45%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.85"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +0.54%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "151.03"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  +0.99%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.841"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +0.35%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.62"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.65%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "285.60"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E44").Value = "  +2.93%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "19.20"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +3.19%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0944"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.16%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0226"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +1.49%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.961.27"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +1.37%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "18.44"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.06%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "111.40"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.20%  "
